$d = $word.ActiveDocument

$quoteOpen = [char]0x201C
$quoteClose = [char]0x201D

# -----------------------------------------------------------------
# Edit 1: "Tweet structure" paragraph - replace the closing sentence
# "Below is an example of Tweet JSON." with new text about the 31
# tweet key-value pairs / 12 classes.
# -----------------------------------------------------------------
$r1 = $d.Content
$newText1 = "The 31 tweet key-value pairs belong to 12 distinct classes (Appendix 1). The classes are either vectors - numeric, logical, or character - or arrays assembled from the vector classes."
$found1 = $r1.Find.Execute("Below is an example of Tweet JSON.", $false, $false, $false, $false, $false, $true, 1, $false, $newText1, 2)
if (-not $found1) { throw "Edit 1: target text not found" }

# -----------------------------------------------------------------
# Edit 2: the "\begin{figure}" paragraph becomes the new
# "Below is an example of Tweet JSON. ..." paragraph.
# -----------------------------------------------------------------
$r2 = $d.Content
$newText2 = "Below is an example of Tweet JSON. Every tweet features the keys " + $quoteOpen + "created_at" + $quoteClose + " (the time stamp), " + $quoteOpen + "id_str" + $quoteClose + " (a unique tweet identifier), and " + $quoteOpen + "text" + $quoteClose + ". We use these three keys in our analyses."
$found2 = $r2.Find.Execute("\begin{figure}", $false, $false, $false, $false, $false, $true, 1, $false, $newText2, 2)
if (-not $found2) { throw "Edit 2: target text not found" }

# -----------------------------------------------------------------
# Edit 3: remove the "\end{figure}" paragraph entirely (its whole
# range, including the paragraph mark), which has the effect of
# moving the "Parsing text of tweets" Heading 3 paragraph directly
# after the JSON source-code block.
# -----------------------------------------------------------------
$r3 = $d.Content
$found3 = $r3.Find.Execute("\end{figure}")
if (-not $found3) { throw "Edit 3: target text not found" }
$p3 = $r3.Paragraphs(1)
$p3.Range.Delete()

# -----------------------------------------------------------------
# Edit 4: rework the "We used functions from the rtweet R package..."
# paragraph, and split a new paragraph off of it describing the
# tidytext word-splitting step.
# -----------------------------------------------------------------
$r4 = $d.Content
$found4 = $r4.Find.Execute("We used functions from the")
if (-not $found4) { throw "Edit 4: target text not found" }
$p4 = $r4.Paragraphs(1)
$p4start = $p4.Range.Start
$p4end = $p4.Range.End

$prefixLen = ("We used functions from the rtweet R package to parse tweet JSON into a data frame. ").Length
$delStart = $p4start + $prefixLen
$delEnd = $p4end - 1
if ($d.Range($delStart, $delEnd).Text -ne "From there, we used tidytext R package functions to break the tweet text into individual words. We discarded commonly used " + $quoteOpen + "stop words" + $quoteClose + " and emojis.") {
    throw "Edit 4: unexpected paragraph content"
}
$d.Range($delStart, $delEnd).Delete()

function Insert-PlainText($pos, $text) {
    $r = $d.Range($pos, $pos)
    $r.InsertAfter($text)
    return $pos + $text.Length
}

function Insert-StyledText($pos, $text, $styleName) {
    $r = $d.Range($pos, $pos)
    $r.InsertAfter($text)
    $r2 = $d.Range($pos, $pos + $text.Length)
    $r2.Style = $styleName
    return $pos + $text.Length
}

$pos = $delStart
$pos = Insert-PlainText $pos "With "
$pos = Insert-StyledText $pos "rtweet" "Verbatim Char"
$pos = Insert-PlainText $pos " functions, we parsed JSON arrays into their component vectors and added them to the data frame."

# split off a new paragraph here
$brRange = $d.Range($pos, $pos)
$brRange.InsertParagraphAfter()
$newParaStart = $pos + 1

$newPara = $d.Range($newParaStart, $newParaStart).Paragraphs(1)
$newPara.Style = "Body Text"

$part1 = "We then divided tweet text into words with functions from the "
$part2 = "tidytext"
$part3 = " R package. We discarded commonly used " + $quoteOpen + "stop words" + $quoteClose + " and emojis."
$fullText = $part1 + $part2 + $part3

$ins = $d.Range($newParaStart, $newParaStart)
$ins.InsertAfter($fullText)

$scopeEnd = $newParaStart + $fullText.Length
$scopedRange = $d.Range($newParaStart, $scopeEnd)
$foundTT = $scopedRange.Find.Execute("tidytext", $true, $true)
if (-not $foundTT) { throw "Edit 4: tidytext substring not found" }
$scopedRange.Style = "Verbatim Char"

$beforeRange = $d.Range($newParaStart, $scopedRange.Start)
$beforeRange.Style = "Default Paragraph Font"
$afterRange = $d.Range($scopedRange.End, $scopeEnd)
$afterRange.Style = "Default Paragraph Font"

# -----------------------------------------------------------------
# Edit 5: "saved it as a csv file" -> "saved it as a supplementary file"
# -----------------------------------------------------------------
$r5 = $d.Content
$found5 = $r5.Find.Execute(", (Accessed: May 23, 2020)). We have saved it as a csv file,", $false, $false, $false, $false, $false, $true, 1, $false, ", (Accessed: May 23, 2020)). We have saved it as a supplementary file,", 2)
if (-not $found5) { throw "Edit 5: target text not found" }
# the replaced run sits right after a hyperlink run and otherwise
# inherits the Hyperlink character style; restore plain formatting.
$r5.Style = "Default Paragraph Font"

Write-Output "All edits applied successfully"
